$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to be stored as text (avoids Excel
# auto-converting numeric-looking strings like "212.72" into real numbers),
# then restores the original "General" number format so styling is unchanged.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

# Row 2
Set-TextValue "D2" "27.951.32"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
Set-TextValue "D3" "1.639.36"
$ws.Range("E3").Value = "  -0.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
Set-TextValue "D5" "212.72"
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("E6").Value = "  -0.18%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
Set-TextValue "D8" "23.37"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9
$ws.Range("E9").Value = "  -2.17%  "

# Row 10
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$ws.Range("E11").Value = "  +1.77%  "

# Row 12
Set-TextValue "D12" "1.872.69"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
Set-TextValue "D13" "1.636.43"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("E14").Value = "  +0.41%  "

# Row 15
Set-TextValue "D15" "0.571"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16
$ws.Range("E16").Value = "  -0.19%  "

# Row 17
Set-TextValue "D17" "27.953.42"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
Set-TextValue "D18" "232.57"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
Set-TextValue "D19" "0.0₃0721"
$ws.Range("E19").Value = "  -0.08%  "

# Row 20
Set-TextValue "D20" "7.57"
$ws.Range("E20").Value = "  -1.05%  "

# Row 21
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
Set-TextValue "D22" "10.52"
$ws.Range("E22").Value = "  -1.96%  "

# Row 23
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
$ws.Range("E24").Value = "  -4.10%  "

# Row 25
Set-TextValue "D25" "153.27"
$ws.Range("E25").Value = "  +1.51%  "

# Row 26
$ws.Range("E26").Value = "  -0.42%  "

# Row 27
$ws.Range("E27").Value = "  -0.32%  "

# Row 28
$ws.Range("E28").Value = "  -0.39%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").Value = "  +0.40%  "

# Row 31
$ws.Range("E31").Value = "  +0.26%  "

# Row 32
$ws.Range("E32").Value = "  +2.64%  "

# Row 33
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D33" "1.406.52"
$ws.Range("E33").Value = "  -4.15%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "3.08"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("E35").Value = "  +1.29%  "

# Row 36
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.44%  "

# Row 38
$ws.Range("E38").Value = "  +0.27%  "

# Row 39
Set-TextValue "D39" "0.878"
$ws.Range("E39").Value = "  -1.34%  "

# Row 40
$ws.Range("E40").Value = "  -0.75%  "

# Row 41
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "67.02"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "1.86"
$ws.Range("E44").Value = "  +3.75%  "

# Row 45
$ws.Range("E45").Value = "  +2.59%  "

# Row 46
$ws.Range("E46").Value = "  -3.21%  "

# Row 47
Set-TextValue "D47" "1.781.09"
$ws.Range("E47").Value = "  -0.26%  "

# Row 48
Set-TextValue "D48" "87.89"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D49" "0.0999"
$ws.Range("E49").Value = "  -0.32%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0506"
$ws.Range("E50").Value = "  +0.06%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.59"
$ws.Range("E51").Value = "  -1.90%  "
